$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.894.16"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "'1.888.52"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'0.7312"
$ws.Range("E5").Value = "  -5.15%  "
$ws.Range("D6").Value = "'242.15"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.3090"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").Value = "'26.28"
$ws.Range("E9").Value = "  -5.72%  "
$ws.Range("D10").Value = "'0.06891"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").Value = "'0.07940"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'0.7636"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").Value = "'1.898.79"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'5.228"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "'91.20"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").Value = "'29.908.56"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "'14.08"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "'5.734"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'239.95"
$ws.Range("E19").Value = "  -6.33%  "
$ws.Range("D20").Value = "'0.000007743"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'2.099.31"
$ws.Range("E22").Value = "  -4.13%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'6.905"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").Value = "'9.285"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("D26").Value = "'164.24"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'18.86"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'0.1267"
$ws.Range("E28").Value = "  -6.08%  "
$ws.Range("E29").Value = "  -11.70%  "
$ws.Range("D30").Value = "'1.356"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").Value = "'1.530"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "'4.291"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").Value = "'4.065"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "'0.05066"
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("D35").Value = "'1.270"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'0.7329"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").Value = "'2.722"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").Value = "'0.01916"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "'6.318"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'74.36"
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("D42").Value = "'0.4426"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'1.924"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'0.8362"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'7.582"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "'100.66"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'9.758"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'37.14"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'2.008.42"
$ws.Range("E50").Value = "  -3.71%  "
$ws.Range("D51").Value = "'942.46"
$ws.Range("E51").Value = "  -3.99%  "
